$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest_eval")

$data = @{
    2 = @(-0.1341439519111721, 0.5935102738008631, 0.6921245334844653, 0.8319402223023391, 0.8463232565061951, 17)
    3 = @(-0.2124529600720004, 0.6536332604753219, 0.6702315204652148, 0.8186766886049797, 0.8149624003509157, 17)
    4 = @(-0.04676944887512436, 0.4652140835355437, 0.3214984223373458, 0.5670083088785788, 0.5824671501905423, 17)
    5 = @(-0.09716908889197996, 0.6607097935511106, 0.7106863686113114, 0.8430221637722886, 0.8631757181035996, 17)
    6 = @(-0.1115222062452062, 0.3482313290500288, 0.1948460206860828, 0.4414136616441349, 0.4402378554120513, 17)
    7 = @(-0.006406154181615844, 0.5272075584481162, 0.4114983443300604, 0.6414813670949925, 0.6611908855342594, 17)
    8 = @(-0.08290538230610117, 0.3521295447676658, 0.1561676696427139, 0.3951805532193023, 0.3982778743278431, 17)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$wb.Save()
